$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 214, shifting existing rows 214:320 down to 215:321
$ws.Rows.Item(214).Insert()

# Populate new row 214 with the new record (matches neighboring Apio entries' formatting)
$ws.Cells.Item(214, 1).Value = 7
$ws.Cells.Item(214, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(214, 3).Value = "Ñuble"
$ws.Cells.Item(214, 4).Value = 45029
$ws.Cells.Item(214, 5).Value = 16
$ws.Cells.Item(214, 6).Value = 100112017
$ws.Cells.Item(214, 7).Value = "Apio"
$ws.Cells.Item(214, 8).Value = "Americana (o)"
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 50
$ws.Cells.Item(214, 11).Value = 9000
$ws.Cells.Item(214, 12).Value = 9000
$ws.Cells.Item(214, 13).Value = 9000
$ws.Cells.Item(214, 14).Value = "`$/docena de matas"
$ws.Cells.Item(214, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(214, 16).Value = 1500
$ws.Cells.Item(214, 17).Value = 6
$ws.Cells.Item(214, 18).Value = "Hortaliza"
